# Update the Btc-Erbb3 LR-pairs sheet with refreshed TPM-derived values.
# A new "ECs" target cluster is introduced, and all NATMI metric columns
# (E:T) are recalculated; rows are re-emitted sorted by Sending/Target cluster.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @"
FAPs,Btc,Erbb3,ECs,1,0.3333333333333333,0.06428666666666667,0.19286,0.005899579586891775,0.008823342375055644,1,0.5,0.08432099999999999,0.168642,0.008155403471034732,0.007703748342576289,0.00542071602,0.03252429612,4.811345184058284E-05,6.797280919781805E-05
FAPs,Btc,Erbb3,FAPs,1,0.3333333333333333,0.06428666666666667,0.19286,0.005899579586891775,0.008823342375055644,3,1,0.1605563333333333,0.481669,0.01552877311895348,0.02200315912062463,0.01032163148222222,0.09289468334000001,9.161323290205167E-05,0.0001941414062540994
FAPs,Btc,Erbb3,Inflammatory-Mac,1,0.3333333333333333,0.06428666666666667,0.19286,0.005899579586891775,0.008823342375055644,3,1,0.5757656666666666,1.727297,0.05568721097278211,0.07890478884789671,0.03701405549111111,0.33312649942,0.000328531133105961,0.0006962039670364651
FAPs,Btc,Erbb3,MuSCs,1,0.3333333333333333,0.06428666666666667,0.19286,0.005899579586891775,0.008823342375055644,2,1,9.042617,18.085234,0.8745886560766319,0.8261529835545378,0.5813197048733334,3.48791822924,0.005159705382316809,0.007289430628075402
FAPs,Btc,Erbb3,Neutrophils,1,0.3333333333333333,0.06428666666666667,0.19286,0.005899579586891775,0.008823342375055644,1,0.3333333333333333,0.07294366666666667,0.218831,0.007055004474844155,0.009996435962300685,0.004689305184444445,0.04220374666,4.16215603852207E-05,8.820197702569778E-05
FAPs,Btc,Erbb3,Resolving-Mac,1,0.3333333333333333,0.06428666666666667,0.19286,0.005899579586891775,0.008823342375055644,3,1,0.4030763333333334,1.209229,0.03898495188575349,0.05523888417206382,0.02591243388222223,0.23321190494,0.0002299948263411493,0.0004873915874661612
MuSCs,Btc,Erbb3,ECs,2,1,10.832535,21.66507,0.9941004204131083,0.9911766576249443,1,0.5,0.08432099999999999,0.168642,0.008155403471034732,0.007703748342576289,0.9134101837349999,3.65364073494,0.008107290019194149,0.007635775533378471
MuSCs,Btc,Erbb3,FAPs,2,1,10.832535,21.66507,0.9941004204131083,0.9911766576249443,3,1,0.1605563333333333,0.481669,0.01552877311895348,0.02200315912062463,1.739232100305,10.43539260183,0.01543715988605143,0.02180901771437053
MuSCs,Btc,Erbb3,Inflammatory-Mac,2,1,10.832535,21.66507,0.9941004204131083,0.9911766576249443,3,1,0.5757656666666666,1.727297,0.05568721097278211,0.07890478884789671,6.237001735964999,37.42201041579,0.05535867983967615,0.07820858488086024
MuSCs,Btc,Erbb3,MuSCs,2,1,10.832535,21.66507,0.9941004204131083,0.9911766576249443,2,1,9.042617,18.085234,0.8745886560766319,0.8261529835545378,97.954465144095,391.81786057638,0.8694289506943153,0.8188635529264624
MuSCs,Btc,Erbb3,Neutrophils,2,1,10.832535,21.66507,0.9941004204131083,0.9911766576249443,1,0.3333333333333333,0.07294366666666667,0.218831,0.007055004474844155,0.009996435962300685,0.7901648221950001,4.74098893317,0.007013382914458934,0.009908233985274987
MuSCs,Btc,Erbb3,Resolving-Mac,2,1,10.832535,21.66507,0.9941004204131083,0.9911766576249443,3,1,0.4030763333333334,1.209229,0.03898495188575349,0.05523888417206382,4.366338488505001,26.19803093103,0.03875495705941234,0.05475149258459765
"@

$lines = $data -split "`n"
$r = 2
foreach ($line in $lines) {
  $line = $line.Trim()
  if ($line.Length -eq 0) { continue }
  $parts = $line -split ","

  $ws.Cells.Item($r, 1).Value = $parts[0]
  $ws.Cells.Item($r, 2).Value = $parts[1]
  $ws.Cells.Item($r, 3).Value = $parts[2]
  $ws.Cells.Item($r, 4).Value = $parts[3]

  for ($c = 5; $c -le 20; $c++) {
    $ws.Cells.Item($r, $c).Value = [double]$parts[$c - 1]
  }

  $r = $r + 1
}
